$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.041269294343564
$ws.Range("D2").Value = 1.051451163315476
$ws.Range("E2").Value = 1.049653614489217
$ws.Range("F2").Value = 1.060723645518419
$ws.Range("I2").Value = 1.039365580196626
$ws.Range("J2").Value = 1.046351413287061
$ws.Range("K2").Value = 1.054202629194361
$ws.Range("L2").Value = 1.052410073043047
$ws.Range("M2").Value = 1.063449648228745
$ws.Range("N2").Value = 1.019361217592652

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.042174561338871
$ws.Range("D3").Value = 1.052288117105676
$ws.Range("E3").Value = 1.050455776351999
$ws.Range("F3").Value = 1.061603838558651
$ws.Range("I3").Value = 1.039510833925993
$ws.Range("J3").Value = 1.046902859859393
$ws.Range("K3").Value = 1.054852169320907
$ws.Range("L3").Value = 1.053024550343679
$ws.Range("M3").Value = 1.06414415671204
$ws.Range("N3").Value = 1.01954421352976

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042760603703636
$ws.Range("D4").Value = 1.052830218481948
$ws.Range("E4").Value = 1.050975432162712
$ws.Range("F4").Value = 1.062173966657135
$ws.Range("I4").Value = 1.039602950373283
$ws.Range("J4").Value = 1.047259312486932
$ws.Range("K4").Value = 1.055272353343205
$ws.Range("L4").Value = 1.053422106243716
$ws.Range("M4").Value = 1.064593493604727
$ws.Range("N4").Value = 1.019662479959742

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043007040221012
$ws.Range("D5").Value = 1.053058244668542
$ws.Range("E5").Value = 1.051194038321909
$ws.Range("F5").Value = 1.062413786712957
$ws.Range("I5").Value = 1.039641227042404
$ws.Range("J5").Value = 1.047409075359062
$ws.Range("K5").Value = 1.055448970567961
$ws.Range("L5").Value = 1.053589224997782
$ws.Range("M5").Value = 1.064782379843548
$ws.Range("N5").Value = 1.019712164224732

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.043048421729569
$ws.Range("D6").Value = 1.053096538655024
$ws.Range("E6").Value = 1.051230751588217
$ws.Range("F6").Value = 1.062454061610661
$ws.Range("I6").Value = 1.039647627504995
$ws.Range("J6").Value = 1.047434215902928
$ws.Range("K6").Value = 1.055478623687612
$ws.Range("L6").Value = 1.053617284128071
$ws.Range("M6").Value = 1.06481409375269
$ws.Range("N6").Value = 1.019720504370012

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042763896349404
$ws.Range("D7").Value = 1.052833264882953
$ws.Range("E7").Value = 1.050978352629616
$ws.Range("F7").Value = 1.062177170603547
$ws.Range("I7").Value = 1.039603463593647
$ws.Range("J7").Value = 1.047261313980164
$ws.Range("K7").Value = 1.055274713423482
$ws.Range("L7").Value = 1.053424339348338
$ws.Range("M7").Value = 1.064596017572638
$ws.Range("N7").Value = 1.019663143981159

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041575176025297
$ws.Range("D8").Value = 1.051733903745588
$ws.Range("E8").Value = 1.049924582932467
$ws.Range("F8").Value = 1.061020988870436
$ws.Range("I8").Value = 1.039415056580492
$ws.Range("J8").Value = 1.046537853171395
$ws.Range("K8").Value = 1.054422166671761
$ws.Range("L8").Value = 1.052617748512885
$ws.Range("M8").Value = 1.06368437145481
$ws.Range("N8").Value = 1.019423091523932

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.039482655699022
$ws.Range("D9").Value = 1.049800867969932
$ws.Range("E9").Value = 1.048072398634682
$ws.Range("F9").Value = 1.058988198696236
$ws.Range("I9").Value = 1.039068754812361
$ws.Range("J9").Value = 1.045260240006682
$ws.Range("K9").Value = 1.052919073163562
$ws.Range("L9").Value = 1.051196091165895
$ws.Range("M9").Value = 1.062077565544311
$ws.Range("N9").Value = 1.018999002907313

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.038089168571993
$ws.Range("D10").Value = 1.048515082745616
$ws.Range("E10").Value = 1.046840858890786
$ws.Range("F10").Value = 1.057636163610659
$ws.Range("I10").Value = 1.038828317332607
$ws.Range("J10").Value = 1.04440669666708
$ws.Range("K10").Value = 1.051916554352355
$ws.Range("L10").Value = 1.050248166308679
$ws.Range("M10").Value = 1.061006198981366
$ws.Range("N10").Value = 1.018715572298724

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.037486150989429
$ws.Range("D11").Value = 1.047959032534492
$ws.Range("E11").Value = 1.046308379155313
$ws.Range("F11").Value = 1.057051487290507
$ws.Range("I11").Value = 1.038721945352506
$ws.Range("J11").Value = 1.044036690205971
$ws.Range("K11").Value = 1.051482360878434
$ws.Range("L11").Value = 1.04983768384283
$ws.Range("M11").Value = 1.060542264611636
$ws.Range("N11").Value = 1.018592681718697

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037262220444253
$ws.Range("D12").Value = 1.047752597965708
$ws.Range("E12").Value = 1.046110711999133
$ws.Range("F12").Value = 1.056834428864605
$ws.Range("I12").Value = 1.038682095133786
$ws.Range("J12").Value = 1.043899191880312
$ws.Range("K12").Value = 1.051321068934477
$ws.Range("L12").Value = 1.049685209908324
$ws.Range("M12").Value = 1.060369936286964
$ws.Range("N12").Value = 1.01854701056582

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037310251705925
$ws.Range("D13").Value = 1.047796874003441
$ws.Range("E13").Value = 1.046153106841099
$ws.Range("F13").Value = 1.056880983340502
$ws.Range("I13").Value = 1.038690658466151
$ws.Range("J13").Value = 1.043928688500583
$ws.Range("K13").Value = 1.051355667178455
$ws.Range("L13").Value = 1.049717916176634
$ws.Range("M13").Value = 1.060406901385716
$ws.Range("N13").Value = 1.018556808268086

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.037467639650203
$ws.Range("D14").Value = 1.047941966391908
$ws.Range("E14").Value = 1.046292037479138
$ws.Range("F14").Value = 1.057033542789475
$ws.Range("I14").Value = 1.038718658230693
$ws.Range("J14").Value = 1.044025325798535
$ws.Range("K14").Value = 1.051469028696893
$ws.Range("L14").Value = 1.049825080343022
$ws.Range("M14").Value = 1.06052801994032
$ws.Range("N14").Value = 1.018588907016636

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037564619151059
$ws.Range("D15").Value = 1.048031376858276
$ws.Range("E15").Value = 1.046377653115899
$ws.Range("F15").Value = 1.057127555173165
$ws.Range("I15").Value = 1.03873586492086
$ws.Range("J15").Value = 1.044084859118127
$ws.Range("K15").Value = 1.051538872765406
$ws.Range("L15").Value = 1.049891107432765
$ws.Range("M15").Value = 1.060602644797286
$ws.Range("N15").Value = 1.018608680931077

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.038129196794574
$ws.Range("D16").Value = 1.048552000918885
$ws.Range("E16").Value = 1.046876214492762
$ws.Range("F16").Value = 1.057674982884862
$ws.Range("I16").Value = 1.038835329319242
$ws.Range("J16").Value = 1.044431244097568
$ws.Range("K16").Value = 1.051945368427172
$ws.Range("L16").Value = 1.050275408263274
$ws.Range("M16").Value = 1.06103698835824
$ws.Range("N16").Value = 1.018723724737

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038483441763931
$ws.Range("D17").Value = 1.04887876404021
$ws.Range("E17").Value = 1.04718916017893
$ws.Range("F17").Value = 1.058018575474707
$ws.Range("I17").Value = 1.038897115918697
$ws.Range("J17").Value = 1.044648411518831
$ws.Range("K17").Value = 1.052200327576229
$ws.Range("L17").Value = 1.050516464277986
$ws.Range("M17").Value = 1.061309434740498
$ws.Range("N17").Value = 1.018795845256987

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03869010262795
$ws.Range("D18").Value = 1.049069427190833
$ws.Range("E18").Value = 1.047371771748617
$ws.Range("F18").Value = 1.058219060804535
$ws.Range("I18").Value = 1.038932936796109
$ws.Range("J18").Value = 1.04477504133352
$ws.Range("K18").Value = 1.052349031564238
$ws.Range("L18").Value = 1.050657065701808
$ws.Range("M18").Value = 1.061468345544334
$ws.Range("N18").Value = 1.018837896157011

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.038760574671667
$ws.Range("D19").Value = 1.049134449841554
$ws.Range("E19").Value = 1.047434050357619
$ws.Range("F19").Value = 1.058287433552353
$ws.Range("I19").Value = 1.03894511375594
$ws.Range("J19").Value = 1.044818211963267
$ws.Range("K19").Value = 1.05239973415236
$ws.Range("L19").Value = 1.050705006696669
$ws.Range("M19").Value = 1.061522529544686
$ws.Range("N19").Value = 1.018852231738393

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.038445430922919
$ws.Range("D20").Value = 1.048843698430579
$ws.Range("E20").Value = 1.047155576232724
$ws.Range("F20").Value = 1.05798170361839
$ws.Range("I20").Value = 1.03889050936536
$ws.Range("J20").Value = 1.044625115666676
$ws.Range("K20").Value = 1.052172973832218
$ws.Range("L20").Value = 1.050490601495506
$ws.Range("M20").Value = 1.061280204076895
$ws.Range("N20").Value = 1.018788109036886

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037421291239226
$ws.Range("D21").Value = 1.047899237326897
$ws.Range("E21").Value = 1.046251122556868
$ws.Range("F21").Value = 1.056988614608052
$ws.Range("I21").Value = 1.03871042234797
$ws.Range("J21").Value = 1.043996870210197
$ws.Range("K21").Value = 1.051435646907909
$ws.Range("L21").Value = 1.049793523225198
$ws.Range("M21").Value = 1.060492353601012
$ws.Range("N21").Value = 1.018579455397786

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036777704556021
$ws.Range("D22").Value = 1.047306038368113
$ws.Range("E22").Value = 1.045683148950668
$ws.Range("F22").Value = 1.056364894349325
$ws.Range("I22").Value = 1.038595233924021
$ws.Range("J22").Value = 1.043601512375277
$ws.Range("K22").Value = 1.050971984734126
$ws.Range("L22").Value = 1.049355228774191
$ws.Range("M22").Value = 1.059996987291863
$ws.Range("N22").Value = 1.018448127224723

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037118850292914
$ws.Range("D23").Value = 1.047620444949281
$ws.Range("E23").Value = 1.045984176357816
$ws.Range("F23").Value = 1.056695475805008
$ws.Range("I23").Value = 1.038656483075181
$ws.Range("J23").Value = 1.043811132315033
$ws.Range("K23").Value = 1.051217787546248
$ws.Range("L23").Value = 1.049587577836083
$ws.Range("M23").Value = 1.060259591153256
$ws.Range("N23").Value = 1.018517759846998

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.038462606275396
$ws.Range("D24").Value = 1.048859542861016
$ws.Range("E24").Value = 1.047170751138846
$ws.Range("F24").Value = 1.057998364195932
$ws.Range("I24").Value = 1.038893495256707
$ws.Range("J24").Value = 1.044635642184035
$ws.Range("K24").Value = 1.052185333840787
$ws.Range("L24").Value = 1.050502287781417
$ws.Range("M24").Value = 1.061293412163611
$ws.Range("N24").Value = 1.018791604750543

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040023358900693
$ws.Range("D25").Value = 1.050300098952236
$ws.Range("E25").Value = 1.048550666795216
$ws.Range("F25").Value = 1.059513174253997
$ws.Range("I25").Value = 1.039159972324514
$ws.Range("J25").Value = 1.045590856204023
$ws.Range("K25").Value = 1.053307745125358
$ws.Range("L25").Value = 1.051563656243487
$ws.Range("M25").Value = 1.062492998287907
$ws.Range("N25").Value = 1.019108766039697
